$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A94").Value = "2023-12-08 09:46:42"
$ws.Range("B94").Value = 0.0004

$ws.Range("A95").Value = "2023-12-08 09:47:27"
$ws.Range("B95").Value = 0.0028

$ws.Range("A96").Value = "2023-12-08 09:48:05"
$ws.Range("B96").Value = 0.003
